# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Vega Monumental Concepción - Limón"
# right above the existing row 288, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 288 (existing rows 288:378 move to 289:379)
$ws.Rows(288).Insert()

# Populate the newly inserted row with the new weekly observation
$ws.Cells.Item(288, 1).Value  = 11
$ws.Cells.Item(288, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(288, 3).Value  = "Bíobío"
$ws.Cells.Item(288, 4).Value  = 44588
$ws.Cells.Item(288, 5).Value  = 8
$ws.Cells.Item(288, 6).Value  = "Fruta"
$ws.Cells.Item(288, 7).Value  = 100102
$ws.Cells.Item(288, 8).Value  = "Cítricos"
$ws.Cells.Item(288, 9).Value  = 100102003
$ws.Cells.Item(288, 10).Value = "Limón"
$ws.Cells.Item(288, 11).Value = "Sin especificar"
$ws.Cells.Item(288, 12).Value = "1a plateado"
$ws.Cells.Item(288, 13).Value = 250
$ws.Cells.Item(288, 14).Value = 14000
$ws.Cells.Item(288, 15).Value = 15000
$ws.Cells.Item(288, 16).Value = 14520
$ws.Cells.Item(288, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(288, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(288, 19).Value = 908
$ws.Cells.Item(288, 20).Value = 16
